$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H holds the "HS" (Highest Score) stat. Values that were recorded as
# text strings with a trailing "*" (denoting a "not out" innings, e.g. "100*")
# are converted to plain numeric values (the "*" designation removed).
$ws.Cells.Item(3, 8).Value = 100
$ws.Cells.Item(4, 8).Value = 2
$ws.Cells.Item(8, 8).Value = 1
$ws.Cells.Item(9, 8).Value = 62
$ws.Cells.Item(13, 8).Value = 2
$ws.Cells.Item(14, 8).Value = 9
$ws.Cells.Item(15, 8).Value = 58
$ws.Cells.Item(17, 8).Value = 76
$ws.Cells.Item(20, 8).Value = 8
$ws.Cells.Item(21, 8).Value = 100
$ws.Cells.Item(25, 8).Value = 1
$ws.Cells.Item(36, 8).Value = 39
$ws.Cells.Item(37, 8).Value = 59
$ws.Cells.Item(39, 8).Value = 68
$ws.Cells.Item(42, 8).Value = 84
$ws.Cells.Item(46, 8).Value = 64
$ws.Cells.Item(48, 8).Value = 25
$ws.Cells.Item(51, 8).Value = 17
$ws.Cells.Item(52, 8).Value = 48
$ws.Cells.Item(59, 8).Value = 59
$ws.Cells.Item(60, 8).Value = 37
$ws.Cells.Item(64, 8).Value = 8
$ws.Cells.Item(66, 8).Value = 10
$ws.Cells.Item(67, 8).Value = 63
$ws.Cells.Item(71, 8).Value = 95
$ws.Cells.Item(72, 8).Value = 13
$ws.Cells.Item(73, 8).Value = 84
$ws.Cells.Item(77, 8).Value = 33
$ws.Cells.Item(78, 8).Value = 75
$ws.Cells.Item(79, 8).Value = 54
$ws.Cells.Item(80, 8).Value = 10
$ws.Cells.Item(81, 8).Value = 20
$ws.Cells.Item(82, 8).Value = 54
$ws.Cells.Item(86, 8).Value = 107
$ws.Cells.Item(87, 8).Value = 91
$ws.Cells.Item(89, 8).Value = 124
$ws.Cells.Item(91, 8).Value = 15
$ws.Cells.Item(98, 8).Value = 12
$ws.Cells.Item(100, 8).Value = 66
$ws.Cells.Item(102, 8).Value = 63
$ws.Cells.Item(104, 8).Value = 25
$ws.Cells.Item(105, 8).Value = 47
$ws.Cells.Item(107, 8).Value = 36
$ws.Cells.Item(108, 8).Value = 88
$ws.Cells.Item(110, 8).Value = 66
$ws.Cells.Item(113, 8).Value = 62
$ws.Cells.Item(115, 8).Value = 106

# Update the active cell selection to reflect where the user left off editing.
$ws.Range("M12").Select()

